$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.120.64"
$ws.Range("E2").Value = "  +2.95%  "
$ws.Range("D3").Value = "3.073.24"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.449"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.68"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.32%  "
$ws.Range("E10").Value = "  +8.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.372"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.32%  "
$ws.Range("E12").Value = "  +2.29%  "
$ws.Range("D13").Value = "3.596.26"
$ws.Range("E13").Value = "  +3.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000172"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +17.42%  "
$ws.Range("D16").Value = "58.096.48"
$ws.Range("E16").Value = "  +2.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.45%  "
$ws.Range("D18").Value = "3.074.57"
$ws.Range("E18").Value = "  +3.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "342.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.504"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  +7.24%  "
$ws.Range("D27").Value = "0.0₃0986"
$ws.Range("E27").Value = "  +10.10%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.38%  "
$ws.Range("E31").Value = "  +7.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.93"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0709"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.87%  "
$ws.Range("D40").Value = "3.108.66"
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.50%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.67%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.668"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.70%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "2.337.82"
$ws.Range("E46").Value = "  +5.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.22%  "
$ws.Range("E48").Value = "  +4.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0247"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.58%  "
